# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" worksheets to reflect the newly scraped totals.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 995
$ws1.Range("F5").Value = 2686
$ws1.Range("F6").Value = 89
$ws1.Range("F7").Value = 204
$ws1.Range("F10").Value = 51
$ws1.Range("F11").Value = 2520
$ws1.Range("F12").Value = 614

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 995
$ws4.Range("F6").Value = 2686
$ws4.Range("F7").Value = 89
$ws4.Range("F8").Value = 204
$ws4.Range("F12").Value = 51
$ws4.Range("F13").Value = 2520
$ws4.Range("F14").Value = 614
